$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F6 date value (optimizacion de las fechas) - keep as text, not a date serial.
# Use a leading apostrophe to force text entry, then clear the resulting
# quote-prefix formatting so no new cell style is introduced.
$ws.Range("F6").Value = "'2024-05-05"
$ws.Range("F6").ClearFormats()

# Add new row 10 with client data
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "ivan"
$ws.Range("C10").Value = "gorda"
$ws.Range("D10").Value = 460036
$ws.Range("E10").Value = "'8675309125"
$ws.Range("E10").ClearFormats()
$ws.Range("F10").Value = "'2024-12-11"
$ws.Range("F10").ClearFormats()
